$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 19 de Julio de 2020 a las 01:34'
$ws.Range("B4").Value = 3831496
$ws.Range("C4").Value = 61484
$ws.Range("D4").Value = 1773783
$ws.Range("E4").Value = 1914853
$ws.Range("G4").Value = 796
$ws.Range("H4").Value = 142860
$ws.Range("B21").Value = 190700
$ws.Range("C21").Value = 8560
$ws.Range("D21").Value = 85836
$ws.Range("E21").Value = 98348
$ws.Range("G21").Value = 228
$ws.Range("H21").Value = 6516
$ws.Range("B31").Value = 73382
$ws.Range("C31").Value = 938
$ws.Range("D31").Value = 31757
$ws.Range("E31").Value = 36343
$ws.Range("G31").Value = 32
$ws.Range("H31").Value = 5282
$ws.Range("A49").Value = 'Nigeria'
$ws.Range("B49").Value = 36107
$ws.Range("C49").Value = 653
$ws.Range("D49").Value = 14938
$ws.Range("E49").Value = 20391
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 778
$ws.Range("A50").Value = 'Barein'
$ws.Range("B50").Value = 36004
$ws.Range("C50").Value = 531
$ws.Range("D50").Value = 31765
$ws.Range("E50").Value = 4115
$ws.Range("H50").Value = 124
$ws.Range("B57").Value = 27060
$ws.Range("C57").Value = 488
$ws.Range("D57").Value = 23044
$ws.Range("E57").Value = 3871
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 145
$ws.Range("A59").Value = 'Kirguistan'
$ws.Range("B59").Value = 24606
$ws.Range("C59").Value = 949
$ws.Range("D59").Value = 10704
$ws.Range("E59").Value = 13002
$ws.Range("G59").Value = 35
$ws.Range("H59").Value = 900
$ws.Range("A60").Value = 'Japon'
$ws.Range("B60").Value = 24132
$ws.Range("C60").Value = 659
$ws.Range("D60").Value = 19366
$ws.Range("E60").Value = 3781
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 985
$ws.Range("A61").Value = 'Argelia'
$ws.Range("B61").Value = 22549
$ws.Range("C61").Value = 601
$ws.Range("D61").Value = 15744
$ws.Range("E61").Value = 5737
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 1068
$ws.Range("A62").Value = 'Moldavia'
$ws.Range("B62").Value = 20794
$ws.Range("C62").Value = 300
$ws.Range("D62").Value = 14183
$ws.Range("E62").Value = 5931
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 680
$ws.Range("A63").Value = 'Serbia'
$ws.Range("B63").Value = 20498
$ws.Range("C63").Value = 389
$ws.Range("D63").Value = 14047
$ws.Range("E63").Value = 5990
$ws.Range("G63").Value = 9
$ws.Range("H63").Value = 461
$ws.Range("A64").Value = 'Austria'
$ws.Range("B64").Value = 19573
$ws.Range("C64").Value = 134
$ws.Range("D64").Value = 17501
$ws.Range("E64").Value = 1361
$ws.Range("H64").Value = 711
$ws.Range("A65").Value = 'Nepal'
$ws.Range("B65").Value = 17502
$ws.Range("C65").Value = 57
$ws.Range("D65").Value = 11637
$ws.Range("E65").Value = 5825
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 40
$ws.Range("A66").Value = 'Marruecos'
$ws.Range("B66").Value = 17015
$ws.Range("C66").Value = 289
$ws.Range("D66").Value = 14620
$ws.Range("E66").Value = 2126
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 269
$ws.Range("A67").Value = 'Uzbekistan'
$ws.Range("B67").Value = 16186
$ws.Range("C67").Value = 579
$ws.Range("D67").Value = 9127
$ws.Range("E67").Value = 6976
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 83
$ws.Range("A68").Value = 'Camerun'
$ws.Range("B68").Value = 16157
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 13728
$ws.Range("E68").Value = 2056
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 373
$ws.Range("A69").Value = 'Costa de Marfil'
$ws.Range("B69").Value = 13912
$ws.Range("C69").Value = 216
$ws.Range("D69").Value = 8000
$ws.Range("E69").Value = 5821
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 91
$ws.Range("A70").Value = 'Chequia'
$ws.Range("B70").Value = 13795
$ws.Range("C70").Value = 53
$ws.Range("D70").Value = 8725
$ws.Range("E70").Value = 4712
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 358
$ws.Range("A71").Value = 'Corea del Sur'
$ws.Range("B71").Value = 13711
$ws.Range("C71").Value = 39
$ws.Range("D71").Value = 12519
$ws.Range("E71").Value = 898
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 294
$ws.Range("A72").Value = 'Dinamarca'
$ws.Range("B72").Value = 13173
$ws.Range("D72").Value = 12209
$ws.Range("E72").Value = 353
$ws.Range("H72").Value = 611
$ws.Range("B93").Value = 6491
$ws.Range("C93").Value = 61
$ws.Range("D93").Value = 5257
$ws.Range("E93").Value = 1195
$ws.Range("B103").Value = 3629
$ws.Range("C103").Value = 172
$ws.Range("D103").Value = 1643
$ws.Range("E103").Value = 1957
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 29
$ws.Range("B106").Value = 3111
$ws.Range("C106").Value = 5
$ws.Range("D106").Value = 1452
$ws.Range("E106").Value = 1566
$ws.Range("A144").Value = 'Uruguay'
$ws.Range("B144").Value = 1044
$ws.Range("C144").Value = 7
$ws.Range("D144").Value = 921
$ws.Range("E144").Value = 90
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 33
$ws.Range("A145").Value = 'Republica de Chipre'
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 845
$ws.Range("E145").Value = 173
$ws.Range("H145").Value = 19
$ws.Range("A180").Value = 'Bahamas'
$ws.Range("B180").Value = 138
$ws.Range("C180").Value = 9
$ws.Range("D180").Value = 91
$ws.Range("E180").Value = 36
$ws.Range("H180").Value = 11
$ws.Range("A181").Value = 'Trinidad yTobago'
$ws.Range("B181").Value = 136
$ws.Range("D181").Value = 124
$ws.Range("E181").Value = 4
$ws.Range("H181").Value = 8
